$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.715.92"
$ws.Range("E2").Value = "  +3.52%  "

$ws.Range("D3").Value = "1.879.44"
$ws.Range("E3").Value = "  +1.78%  "

$ws.Range("D4").Value = "'1.021"
$ws.Range("E4").Value = "  +1.24%  "

$ws.Range("D5").Value = "'336.68"
$ws.Range("E5").Value = "  +0.44%  "

$ws.Range("D6").Value = "'1.015"
$ws.Range("E6").Value = "  +0.82%  "

$ws.Range("D7").Value = "'0.4593"
$ws.Range("E7").Value = "  -1.61%  "

$ws.Range("D8").Value = "'0.3968"
$ws.Range("E8").Value = "  +3.11%  "

$ws.Range("D9").Value = "'48.25"
$ws.Range("E9").Value = "  +2.98%  "

$ws.Range("D10").Value = "'0.07957"
$ws.Range("E10").Value = "  +0.65%  "

$ws.Range("D11").Value = "'0.9933"
$ws.Range("E11").Value = "  +2.84%  "

$ws.Range("D12").Value = "'21.75"
$ws.Range("E12").Value = "  +2.29%  "

$ws.Range("D13").Value = "1.925.46"
$ws.Range("E13").Value = "  +2.78%  "

$ws.Range("D14").Value = "'5.930"
$ws.Range("E14").Value = "  +1.18%  "

$ws.Range("D15").Value = "'7.091"
$ws.Range("E15").Value = "  -0.36%  "

$ws.Range("D16").Value = "'1.027"
$ws.Range("E16").Value = "  +1.87%  "

$ws.Range("D17").Value = "'89.04"
$ws.Range("E17").Value = "  -2.03%  "

$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "'0.06634"
$ws.Range("E18").Value = "  +0.60%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.00001035"
$ws.Range("E19").Value = "  +0.71%  "

$ws.Range("D20").Value = "'17.29"
$ws.Range("E20").Value = "  +0.48%  "

$ws.Range("D21").Value = "'1.012"
$ws.Range("E21").Value = "  +0.51%  "

$ws.Range("D22").Value = "28.844.29"
$ws.Range("E22").Value = "  +3.95%  "

$ws.Range("D23").Value = "'5.431"
$ws.Range("E23").Value = "  +1.79%  "

$ws.Range("D24").Value = "'10.85"
$ws.Range("E24").Value = "  +0.51%  "

$ws.Range("D25").Value = "'2.272"
$ws.Range("E25").Value = "  -1.07%  "

$ws.Range("D26").Value = "2.145.45"
$ws.Range("E26").Value = "  +2.53%  "

$ws.Range("D27").Value = "'159.28"
$ws.Range("E27").Value = "  +0.13%  "

$ws.Range("D28").Value = "'19.59"
$ws.Range("E28").Value = "  +0.86%  "

$ws.Range("D29").Value = "'2.086"
$ws.Range("E29").Value = "  +1.21%  "

$ws.Range("D30").Value = "'5.372"
$ws.Range("E30").Value = "  +0.26%  "

$ws.Range("D31").Value = "'119.12"
$ws.Range("E31").Value = "  +0.56%  "

$ws.Range("D32").Value = "'0.9637"
$ws.Range("E32").Value = "  +2.85%  "

$ws.Range("D33").Value = "'0.09464"
$ws.Range("E33").Value = "  +0.36%  "

$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'3.621"
$ws.Range("E34").Value = "  +0.59%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.414"
$ws.Range("E35").Value = "  +6.76%  "

$ws.Range("D36").Value = "'5.313"
$ws.Range("E36").Value = "  +1.27%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.06081"
$ws.Range("E37").Value = "  +1.07%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02238"
$ws.Range("E38").Value = "  +1.60%  "

$ws.Range("D39").Value = "'8.221"
$ws.Range("E39").Value = "  +0.21%  "

$ws.Range("D40").Value = "'1.162"
$ws.Range("E40").Value = "  +0.39%  "

$ws.Range("B41").Value = "Frax"
$ws.Range("C41").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D41").Value = "'1.014"
$ws.Range("E41").Value = "  +0.79%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.5831"
$ws.Range("E42").Value = "  +0.61%  "

$ws.Range("D43").Value = "'10.19"
$ws.Range("E43").Value = "  +1.55%  "

$ws.Range("D44").Value = "'0.1832"
$ws.Range("E44").Value = "  -0.57%  "

$ws.Range("D45").Value = "'1.257"
$ws.Range("E45").Value = "  -2.07%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'12.26"
$ws.Range("E46").Value = "  +2.10%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "'2.291"
$ws.Range("E47").Value = "  +28.69%  "

$ws.Range("D48").Value = "'0.5489"
$ws.Range("E48").Value = "  +0.96%  "

$ws.Range("D49").Value = "'0.07272"
$ws.Range("E49").Value = "  +6.53%  "

$ws.Range("D50").Value = "'1.905"
$ws.Range("E50").Value = "  -1.21%  "

$ws.Range("D51").Value = "'110.86"
$ws.Range("E51").Value = "  +0.15%  "
